$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.611.53'
$ws.Range("E2").Value = '  -1.08%  '

$ws.Range("D3").Value = '2.285.59'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.48'
$ws.Range("E5").Value = '  +1.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.35'
$ws.Range("E6").Value = '  -3.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -3.11%  '

$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  -3.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.93'
$ws.Range("E10").Value = '  -3.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0784'
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.28'
$ws.Range("E12").Value = '  +3.10%  '

$ws.Range("E13").Value = '  +1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.68'
$ws.Range("E14").Value = '  -2.71%  '

$ws.Range("D15").Value = '2.643.69'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").Value = '2.296.37'
$ws.Range("E16").Value = '  -2.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.775'
$ws.Range("E17").Value = '  -1.69%  '

$ws.Range("D18").Value = '42.479.87'
$ws.Range("E18").Value = '  -1.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.84'
$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  -2.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.00'
$ws.Range("E21").Value = '  -1.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.83'
$ws.Range("E22").Value = '  -3.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.51'
$ws.Range("E23").Value = '  -1.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  -0.88%  '

$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.88'
$ws.Range("E28").Value = '  -0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.56'
$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.04'
$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.98'
$ws.Range("E31").Value = '  -1.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.55'
$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  -1.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.93'
$ws.Range("E35").Value = '  -2.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.44'
$ws.Range("E36").Value = '  -2.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("E37").Value = '  -0.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0688'
$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.73'
$ws.Range("E40").Value = '  -2.40%  '

$ws.Range("E41").Value = '  -2.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.67'
$ws.Range("E42").Value = '  -3.79%  '

$ws.Range("D43").Value = '1.990.10'
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0277'
$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.01'
$ws.Range("E45").Value = '  +2.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.99'
$ws.Range("E46").Value = '  -3.09%  '

$ws.Range("E47").Value = '  -9.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.77'
$ws.Range("E48").Value = '  -2.29%  '

$ws.Range("E49").Value = '  +8.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.40'
$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").Value = '2.510.38'
$ws.Range("E51").Value = '  -0.53%  '
